# Insert a new data row before the current row 389 (shifts 389-413 down to
# 390-414) and populate the new row with the new weekly Apio price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(389).Insert()

$ws.Cells.Item(389,1).Value = 3
$ws.Cells.Item(389,2).Value = "Femacal de La Calera"
$ws.Cells.Item(389,3).Value = "Coquimbo"
$ws.Cells.Item(389,4).Value = 44746
$ws.Cells.Item(389,5).Value = 5
$ws.Cells.Item(389,6).Value = 100112017
$ws.Cells.Item(389,7).Value = "Apio"
$ws.Cells.Item(389,8).Value = "Americana (o)"
$ws.Cells.Item(389,9).Value = "Primera"
$ws.Cells.Item(389,10).Value = 323
$ws.Cells.Item(389,11).Value = 9000
$ws.Cells.Item(389,12).Value = 10000
$ws.Cells.Item(389,13).Value = 9503
$ws.Cells.Item(389,14).Value = '$/docena de matas'
$ws.Cells.Item(389,15).Value = "Pan de Azúcar"
$ws.Cells.Item(389,16).Value = 1584
$ws.Cells.Item(389,17).Value = 6
$ws.Cells.Item(389,18).Value = "Hortaliza"
